# Adds season-record columns (Wins, Losses, Ties) to the Cleveland 2006 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (AC1) onto the new header cells
# so the new headers get the same bold/centered/bordered formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header labels for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-47): 78 wins, 84 losses, 0 ties
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 84
    $ws.Cells.Item($r, 32).Value = 0
}
